$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.80176279249223
$ws.Range("C2").Value = 11.31560241107106
$ws.Range("D2").Value = 10.32841491369376
$ws.Range("E2").Value = 9.492487960912264
$ws.Range("F2").Value = 10.7197553510972
$ws.Range("G2").Value = 9.863264076340414
$ws.Range("H2").Value = 9.642843714138213
$ws.Range("I2").Value = 9.387057812681991
$ws.Range("J2").Value = 10.01160140298562
$ws.Range("K2").Value = 9.953789848704155
$ws.Range("L2").Value = 11.3758289578403
$ws.Range("M2").Value = 10.4590705843998
$ws.Range("N2").Value = 9.905485443912186
$ws.Range("B3").Value = -0.5872448126787629
$ws.Range("C3").Value = -0.6564743383561781
$ws.Range("D3").Value = -0.5853904730264068
$ws.Range("E3").Value = -0.4632241001598181
$ws.Range("F3").Value = -0.6567409847010769
$ws.Range("G3").Value = -0.5561258941997722
$ws.Range("H3").Value = -0.3519228631559845
$ws.Range("I3").Value = -0.4732543789404098
$ws.Range("J3").Value = -0.5595401934443966
$ws.Range("K3").Value = -0.530017553600579
$ws.Range("L3").Value = -0.5904889784709213
$ws.Range("M3").Value = -0.5898823165375545
$ws.Range("N3").Value = -0.5111154817230391
$ws.Range("B4").Value = 0.01122169160202996
$ws.Range("C4").Value = 0.05693857251431406
$ws.Range("D4").Value = 0.1852896974756091
$ws.Range("E4").Value = 0.1379507311339778
$ws.Range("F4").Value = 0.08819535924512301
$ws.Range("G4").Value = 0.2152312666512526
$ws.Range("H4").Value = 0.1067470568024271
$ws.Range("I4").Value = 0.3149327272050564
$ws.Range("J4").Value = 0.153345788175705
$ws.Range("K4").Value = 0.1774838196556885
$ws.Range("L4").Value = -0.08296074681812737
$ws.Range("M4").Value = 0.05272847779437562
$ws.Range("N4").Value = 0.1958041521628128
$ws.Range("B5").Value = 0.2305985278385058
$ws.Range("C5").Value = 0.4797485940178183
$ws.Range("D5").Value = 0.7868956032513934
$ws.Range("E5").Value = 0.3979449685933885
$ws.Range("F5").Value = 0.6264893266796507
$ws.Range("G5").Value = 0.4686160949836983
$ws.Range("H5").Value = 0.3777420181051386
$ws.Range("I5").Value = 0.7095609614332354
$ws.Range("J5").Value = 0.4978861829254566
$ws.Range("K5").Value = 0.4373657267465099
$ws.Range("L5").Value = 0.5128784138447524
$ws.Range("M5").Value = 0.396092863045816
$ws.Range("N5").Value = 0.6869262110785017
$ws.Range("B6").Value = 0.02274971694723556
$ws.Range("C6").Value = 0.04732963745277272
$ws.Range("D6").Value = 0.07763125120005265
$ws.Range("E6").Value = 0.03925929398134009
$ws.Range("F6").Value = 0.06180635664078315
$ws.Range("G6").Value = 0.04623135983445693
$ws.Range("H6").Value = 0.03860358192356331
$ws.Range("I6").Value = 0.07000179567806834
$ws.Range("J6").Value = 0.04911900279530892
$ws.Range("K6").Value = 0.04314835215632116
$ws.Range("L6").Value = 0.05059806258384509
$ws.Range("M6").Value = 0.03907657435446898
$ws.Range("N6").Value = 0.06776876249885327
$ws.Range("B7").Value = 0.03873885909775399
$ws.Range("C7").Value = 0.08059424038915525
$ws.Range("D7").Value = 0.1321926821681455
$ws.Range("E7").Value = 0.06685183210621187
$ws.Range("F7").Value = 0.1052456057719775
$ws.Range("G7").Value = 0.07872406231156284
$ws.Range("H7").Value = 0.06734014259899682
$ws.Range("I7").Value = 0.1192010303096091
$ws.Range("J7").Value = 0.0836412221181889
$ws.Range("K7").Value = 0.07347422995902854
$ws.Range("L7").Value = 0.08615980680555659
$ws.Range("M7").Value = 0.06654069197659306
$ws.Range("N7").Value = 0.1153985584858559
$ws.Range("G10").Value = "*"
